$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.008.30'
$ws.Range('E2').Value = '  +0.39%  '

$ws.Range('D3').Value = '2.552.52'
$ws.Range('E3').Value = '  +0.63%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''305.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.08%  '

$ws.Range('D6').Value = '''98.25'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.69%  '

$ws.Range('D7').Value = '''0.579'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.19%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('E9').Value = '  -0.36%  '

$ws.Range('D10').Value = '''37.05'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.51%  '

$ws.Range('D11').Value = '''0.0833'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.72%  '

$ws.Range('D12').Value = '''7.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.12%  '

$ws.Range('E13').Value = '  +2.04%  '

$ws.Range('D14').Value = '2.945.19'
$ws.Range('E14').Value = '  +0.85%  '

$ws.Range('D15').Value = '2.547.17'
$ws.Range('E15').Value = '  +0.99%  '

$ws.Range('D16').Value = '''15.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.13%  '

$ws.Range('D17').Value = '''0.873'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.54%  '

$ws.Range('D18').Value = '43.012.67'
$ws.Range('E18').Value = '  +0.38%  '

$ws.Range('D19').Value = '''13.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.79%  '

$ws.Range('D20').Value = '0.0₃0997'
$ws.Range('E20').Value = '  +1.49%  '

$ws.Range('D21').Value = '''6.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.59%  '

$ws.Range('D22').Value = '''72.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.78%  '

$ws.Range('D23').Value = '''255.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.37%  '

$ws.Range('D24').Value = '''2.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.75%  '

$ws.Range('D25').Value = '''2.10'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.41%  '

$ws.Range('D26').Value = '''28.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.87%  '

$ws.Range('E27').Value = '  -0.13%  '

$ws.Range('D28').Value = '''10.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.63%  '

$ws.Range('D29').Value = '''37.96'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.16%  '

$ws.Range('D30').Value = '''6.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.35%  '

$ws.Range('D31').Value = '''2.09'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.43%  '

$ws.Range('D32').Value = '''158.60'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.09%  '

$ws.Range('D33').Value = '''19.68'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +16.83%  '

$ws.Range('E34').Value = '  -0.20%  '

$ws.Range('E35').Value = '  +1.23%  '

$ws.Range('D36').Value = '''3.32'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.70%  '

$ws.Range('D37').Value = '''2.63'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.31%  '

$ws.Range('E38').Value = '  +2.97%  '

$ws.Range('D39').Value = '''25.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.36%  '

$ws.Range('E40').Value = '  +0.26%  '

$ws.Range('E41').Value = '  +1.21%  '

$ws.Range('E42').Value = '  +0.83%  '

$ws.Range('E43').Value = '  +28.91%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.097.37'
$ws.Range('E44').Value = '  +0.51%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0307'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.95%  '

$ws.Range('D46').Value = '''0.999'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('D47').Value = '''86.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.64%  '

$ws.Range('D48').Value = '''9.02'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.41%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.802.59'
$ws.Range('E49').Value = '  +0.87%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '''75.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.39%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''103.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.01%  '
